# Generate Report for Handoff
# Adds two new localized-file rows (b66a1510-... and f5dbdca9-...) to each of
# the three report sheets (Overview, zh-cn, de-de), expanding their tables
# from 3 rows (1 header + 2 data) to 5 rows (1 header + 4 data).

$wb = $excel.ActiveWorkbook

$file1Name = "b66a1510-a3f8-4093-8ab6-e26ec4d3222f.md"
$file1Path = "e2e\b66a1510-a3f8-4093-8ab6-e26ec4d3222f.md"
$file2Name = "f5dbdca9-ab15-4a7c-8c58-316181ee9ce7.md"
$file2Path = "e2e\f5dbdca9-ab15-4a7c-8c58-316181ee9ce7.md"

$file1UrlSrc = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e152937fc5b8011e332eeec356cfd0e7fc17273/e2e/b66a1510-a3f8-4093-8ab6-e26ec4d3222f.md"
$file2UrlSrc = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d895a092031391b258f5a9322ecf885cba9e9b/e2e/f5dbdca9-ab15-4a7c-8c58-316181ee9ce7.md"

$genDate = "2016-08-20 18:51:22"

# ---------------------------------------------------------------------------
# Sheet "Overview" -- File Name | Path And Name | Extension | Publish URL |
#                     zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $file1Name
$wsOverview.Range("B4").Value = $file1Path
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = $genDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value = $file2Name
$wsOverview.Range("B5").Value = $file2Path
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = $genDate
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $file1UrlSrc, "", "", $file1Path) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $file2UrlSrc, "", "", $file2Path) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -- Source File Name | File Extension | Status | Source Path |
#                  Priority | Content Duplicate | Latest Handoff File |
#                  Latest Handoff Datetime | Latest Target File |
#                  Latest Handback File | Latest Handback DateTime |
#                  Reference Tokens | To be localized | Dependency From |
#                  Has metadata | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = $file1Name
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "b66a1510-a3f8-4093-8ab6-e26ec4d3222f.4e34c681e7e64bc6d66e6c277c3269fd14f5c24f.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-20 18:51:17"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Range("A5").Value = $file2Name
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "f5dbdca9-ab15-4a7c-8c58-316181ee9ce7.9214fd03ccd40e73970d78d17c2f07b45f0f587a.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-20 18:51:17"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L5").Value = ""
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "False"
$wsZhCn.Range("P5").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $file1UrlSrc, "", "", $file1Name) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), $file2UrlSrc, "", "", $file2Name) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" -- same column layout as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = $file1Name
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "b66a1510-a3f8-4093-8ab6-e26ec4d3222f.4e34c681e7e64bc6d66e6c277c3269fd14f5c24f.de-de.xlf"
$wsDeDe.Range("H4").Value = $genDate
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Range("A5").Value = $file2Name
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "f5dbdca9-ab15-4a7c-8c58-316181ee9ce7.9214fd03ccd40e73970d78d17c2f07b45f0f587a.de-de.xlf"
$wsDeDe.Range("H5").Value = $genDate
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L5").Value = ""
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "False"
$wsDeDe.Range("P5").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $file1UrlSrc, "", "", $file1Name) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), $file2UrlSrc, "", "", $file2Name) | Out-Null
